$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "yShift"
$ws.Range("H1").Value = "sprite height"

$ws.Range("G2").Formula = "=1/100"
$ws.Range("H2").Value = 15

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 15

$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 15

$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 15

$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 19

$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 15

$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 19

$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 15

$ws.Range("H2").Select()
